$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2036.7894
$ws.Range("I32").Value = 1340.2
$ws.Range("J32").Value = 2285.5715
$ws.Range("K32").Value = 1340.2
$ws.Range("L32").Value = 2285.5715
$ws.Range("M32").Value = -1014.2
$ws.Range("N32").Value = -2937.5715
$ws.Range("H43").Value = 1882.4166
$ws.Range("I43").Value = 1200
$ws.Range("J43").Value = 2018.9
$ws.Range("K43").Value = 1200
$ws.Range("L43").Value = 2018.9
$ws.Range("M43").Value = -1131
$ws.Range("N43").Value = -2156.9
$ws.Range("H106").Value = 60608360
$ws.Range("I106").Value = 33335874
$ws.Range("J106").Value = 83335430
$ws.Range("K106").Value = 33335874
$ws.Range("L106").Value = 83335430
$ws.Range("M106").Value = -33335243
$ws.Range("N106").Value = -83336692
$ws.Range("H131").Value = 628.9
$ws.Range("I131").Value = 628.9
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 1886.7
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = 3153.3
$ws.Range("N131").ClearContents()
$ws.Range("H137").Value = 1653.381
$ws.Range("I137").Value = 1501.2354
$ws.Range("J137").Value = 2300
$ws.Range("K137").Value = 4503.706200000001
$ws.Range("L137").Value = 6900
$ws.Range("M137").Value = -1953.706200000001
$ws.Range("N137").Value = -12000
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5481.744
$ws.Range("I32").Value = 3819.6858
$ws.Range("K32").Value = 3819.6858
$ws.Range("M32").Value = -3532.6858
$ws.Range("H122").Value = 1710714.4
$ws.Range("I122").Value = 3206279.5
$ws.Range("K122").Value = 9618838.5
$ws.Range("M122").Value = -9616388.5
$ws.Range("H132").Value = 2341.6038
$ws.Range("I132").Value = 1468.2683
$ws.Range("J132").Value = 5325.5
$ws.Range("K132").Value = 4404.8049
$ws.Range("L132").Value = 15976.5
$ws.Range("M132").Value = -1874.8049
$ws.Range("N132").Value = -21036.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 45455330
$ws.Range("I99").Value = 55556360
$ws.Range("J99").Value = 675.25
$ws.Range("K99").Value = 55556360
$ws.Range("L99").Value = 675.25
$ws.Range("M99").Value = -55554862
$ws.Range("N99").Value = -3671.25
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4050011.2
$ws.Range("I16").Value = 9616370
$ws.Range("J16").Value = 1750
$ws.Range("K16").Value = 9616370
$ws.Range("L16").Value = 1750
$ws.Range("M16").Value = -9616083
$ws.Range("N16").Value = -2324
$ws.Range("H31").Value = 14289515
$ws.Range("I31").Value = 2040.4762
$ws.Range("J31").Value = 35720730
$ws.Range("K31").Value = 2040.4762
$ws.Range("L31").Value = 35720730
$ws.Range("M31").Value = -1745.4762
$ws.Range("N31").Value = -35721320
$ws.Range("H34").Value = 14289515
$ws.Range("I34").Value = 2040.4762
$ws.Range("J34").Value = 35720730
$ws.Range("K34").Value = 2040.4762
$ws.Range("L34").Value = 35720730
$ws.Range("M34").Value = -1838.4762
$ws.Range("N34").Value = -35721134
$ws.Range("H99").Value = 4170066.5
$ws.Range("I99").Value = 2596
$ws.Range("J99").Value = 41677300
$ws.Range("K99").Value = 2596
$ws.Range("L99").Value = 41677300
$ws.Range("M99").Value = -1098
$ws.Range("N99").Value = -41680296
$ws.Range("H113").Value = 4050011.2
$ws.Range("I113").Value = 9616370
$ws.Range("J113").Value = 1750
$ws.Range("K113").Value = 9616370
$ws.Range("L113").Value = 1750
$ws.Range("M113").Value = -9614200
$ws.Range("N113").Value = -6090
$ws.Range("H126").Value = 4170066.5
$ws.Range("I126").Value = 2596
$ws.Range("J126").Value = 41677300
$ws.Range("K126").Value = 7788
$ws.Range("L126").Value = 125031900
$ws.Range("M126").Value = -5318
$ws.Range("N126").Value = -125036840
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 800
$ws.Range("J17").Value = 800
$ws.Range("L17").Value = 2400
$ws.Range("N17").Value = -2738
$ws.Range("H23").Value = 6250093.5
$ws.Range("J23").Value = 97.454544
$ws.Range("L23").Value = 292.363632
$ws.Range("N23").Value = -762.3636320000001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6790.909
$ws.Range("I70").Value = 6400
$ws.Range("K70").Value = 6400
$ws.Range("M70").Value = -6130
$ws.Range("H73").Value = 6790.909
$ws.Range("I73").Value = 6400
$ws.Range("K73").Value = 6400
$ws.Range("M73").Value = -5464
$ws.Range("H102").Value = 1276.3667
$ws.Range("I102").Value = 1136.5652
$ws.Range("K102").Value = 1136.5652
$ws.Range("M102").Value = 485.4348
$ws.Range("H113").Value = 27778832
$ws.Range("H122").Value = 4418505
$ws.Range("I122").Value = 3814854.2
$ws.Range("J122").Value = 5558734
$ws.Range("K122").Value = 11444562.6
$ws.Range("L122").Value = 16676202
$ws.Range("M122").Value = -11442112.6
$ws.Range("N122").Value = -16681102
$ws.Range("H132").Value = 5955657
$ws.Range("I132").Value = 8775302
$ws.Range("J132").Value = 3072.2222
$ws.Range("K132").Value = 26325906
$ws.Range("L132").Value = 9216.6666
$ws.Range("M132").Value = -26323376
$ws.Range("N132").Value = -14276.6666
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2530.5715
$ws.Range("I7").Value = 2379.8
$ws.Range("J7").Value = 2907.5
$ws.Range("K7").Value = 2379.8
$ws.Range("L7").Value = 2907.5
$ws.Range("M7").Value = -2267.8
$ws.Range("N7").Value = -3131.5
$ws.Range("H126").Value = 2530.5715
$ws.Range("I126").Value = 2379.8
$ws.Range("J126").Value = 2907.5
$ws.Range("K126").Value = 7139.400000000001
$ws.Range("L126").Value = 8722.5
$ws.Range("M126").Value = -4669.400000000001
$ws.Range("N126").Value = -13662.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 23608.334
$ws.Range("J124").Value = 23608.334
$ws.Range("L124").Value = 23608.334
$ws.Range("N124").Value = -33428.334
$ws.Range("H132").Value = 1183.9482
$ws.Range("I132").Value = 898
$ws.Range("J132").Value = 1934.5625
$ws.Range("K132").Value = 2694
$ws.Range("L132").Value = 5803.6875
$ws.Range("M132").Value = -164
$ws.Range("N132").Value = -10863.6875
